$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "About"
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

# Fully clear the empty "spacer" cells -- they drop out of the sheet entirely
# (no more <row> element at all for rows 5, 12, 19, 22).
foreach ($r in @(5, 12, 19, 22)) {
    $wsAbout.Range("A$r").Clear()
}

# These bullet cells lost their (cosmetic / no-visual-effect) cell style --
# reset them back to the workbook's default "Normal" style.
foreach ($r in @(9, 10, 11, 14, 15, 16, 17, 18, 21, 23, 24, 25, 26, 27, 28, 29, 30)) {
    $wsAbout.Range("A$r").Style = "Normal"
}

# New "BR 2024 Update" note block at the bottom of the sheet.
$wsAbout.Range("A33").Value = "BR 2024 Update"
$wsAbout.Range("A33").Font.Bold = $true
$wsAbout.Range("A34").Value = "Using [eps-brazil-2.1.1] values in both tables"

# Update the view/selection to match the edited state.
$wsAbout.Range("A35").Select()

# ---------------------------------------------------------------------------
# Sheet "BPaFF-BITPTaP"
# ---------------------------------------------------------------------------
$wsBITPTaP = $wb.Worksheets.Item("BPaFF-BITPTaP")

# Lost (no-visual-effect) cell styles -- reset to "Normal".
$wsBITPTaP.Range("A2").Style = "Normal"
$wsBITPTaP.Range("A11").Style = "Normal"

# These cells used to mirror another cell via a formula; they are now plain
# hard-coded values (cached results kept, formulas removed).
$wsBITPTaP.Range("B13").Value = 1
$wsBITPTaP.Range("B14").Value = 0
$wsBITPTaP.Range("B15").Value = 1
$wsBITPTaP.Range("B16").Value = 1
$wsBITPTaP.Range("B17").Value = 1

$wsBITPTaP.Columns("A").ColumnWidth = 26.45

$wsBITPTaP.Range("B2:B17").Select()

# ---------------------------------------------------------------------------
# Sheet "BPaFF-BDTPTPF"
# ---------------------------------------------------------------------------
$wsBDTPTPF = $wb.Worksheets.Item("BPaFF-BDTPTPF")

$wsBDTPTPF.Range("A2").Style = "Normal"

$wsBDTPTPF.Range("B13").Value = 1
$wsBDTPTPF.Range("B14").Value = 0
$wsBDTPTPF.Range("B15").Value = 1
$wsBDTPTPF.Range("B16").Value = 1
$wsBDTPTPF.Range("B17").Value = 1

$wsBDTPTPF.Columns("A").ColumnWidth = 26.45

$wsBDTPTPF.Range("D18").Select()

# Leave the "About" sheet active/selected, matching tabSelected="1" in sheet1.
$wsAbout.Activate()
